{"js": "// Add month/year to each job's date range in the PROFESSIONAL EXPERIENCE\n// section headings, e.g. \"CGI \u2013 Daugherty / Edward Jones (2021\u20132024)\" ->\n// \"CGI \u2013 Daugherty / Edward Jones (January 2021 \u2013 December 2024)\".\nconst replacements = [\n  [\n    \"CGI \u2013 Daugherty / Edward Jones (2021\u20132024)\",\n    \"CGI \u2013 Daugherty / Edward Jones (January 2021 \u2013 December 2024)\"\n  ],\n  [\n    \"Daugherty \u2013 Cox Communications (2021\u20132024)\",\n    \"Daugherty \u2013 Cox Communications (January 2021 \u2013 December 2024)\"\n  ],\n  [\n    \"BPM Software Solutions (2017\u20132021)\",\n    \"BPM Software Solutions (January 2017 \u2013 December 2021)\"\n  ],\n  [\n    \"Interactive Business Solutions (2016\u20132017)\",\n    \"Interactive Business Solutions (January 2016 \u2013 December 2017)\"\n  ],\n  [\n    \"Soave Enterprises (2015 \u2013 2016)\",\n    \"Soave Enterprises (January 2015 \u2013 December 2016)\"\n  ],\n  [\n    \"John Deere Landscapes (2010 \u2013 2015)\",\n    \"John Deere Landscapes (January 2010 \u2013 December 2015)\"\n  ],\n  [\n    \"Compuware (2006 \u2013 2010)\",\n    \"Compuware (January 2006 \u2013 December 2010)\"\n  ],\n  [\n    \"ParTech, Inc. (2002 \u2013 2005)\",\n    \"ParTech, Inc. (January 2002 \u2013 December 2005)\"\n  ],\n  [\n    \"Nexiq Technologies (formerly MPSI) (1999 \u2013 2002)\",\n    \"Nexiq Technologies (formerly MPSI) (January 1999 \u2013 December 2002)\"\n  ]\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const hit of results.items) {\n    hit.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Add month/year to each job's date range in the PROFESSIONAL EXPERIENCE\n# section headings, e.g. \"CGI \u2013 Daugherty / Edward Jones (2021\u20132024)\" ->\n# \"CGI \u2013 Daugherty / Edward Jones (January 2021 \u2013 December 2024)\".\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"CGI \u2013 Daugherty / Edward Jones (2021\u20132024)\", \"CGI \u2013 Daugherty / Edward Jones (January 2021 \u2013 December 2024)\"),\n    @(\"Daugherty \u2013 Cox Communications (2021\u20132024)\", \"Daugherty \u2013 Cox Communications (January 2021 \u2013 December 2024)\"),\n    @(\"BPM Software Solutions (2017\u20132021)\", \"BPM Software Solutions (January 2017 \u2013 December 2021)\"),\n    @(\"Interactive Business Solutions (2016\u20132017)\", \"Interactive Business Solutions (January 2016 \u2013 December 2017)\"),\n    @(\"Soave Enterprises (2015 \u2013 2016)\", \"Soave Enterprises (January 2015 \u2013 December 2016)\"),\n    @(\"John Deere Landscapes (2010 \u2013 2015)\", \"John Deere Landscapes (January 2010 \u2013 December 2015)\"),\n    @(\"Compuware (2006 \u2013 2010)\", \"Compuware (January 2006 \u2013 December 2010)\"),\n    @(\"ParTech, Inc. (2002 \u2013 2005)\", \"ParTech, Inc. (January 2002 \u2013 December 2005)\"),\n    @(\"Nexiq Technologies (formerly MPSI) (1999 \u2013 2002)\", \"Nexiq Technologies (formerly MPSI) (January 1999 \u2013 December 2002)\")\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    $range.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n"}
